# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "57.691.84";  E = "  +0.32%  " }
    @{ Row = 3;  D = "3.121.61";   E = "  +0.46%  " }
    @{ Row = 4;  D = $null;        E = "  +0.03%  " }
    @{ Row = 5;  D = "532.22";     E = "  +1.19%  " }
    @{ Row = 6;  D = "138.08";     E = "  +0.85%  " }
    @{ Row = 7;  D = $null;        E = "  +0.10%  " }
    @{ Row = 8;  D = "3.119.75";   E = "  +0.51%  " }
    @{ Row = 9;  D = $null;        E = "  +6.60%  " }
    @{ Row = 10; D = $null;        E = "  +0.32%  " }
    @{ Row = 11; D = $null;        E = "  +0.14%  " }
    @{ Row = 12; D = $null;        E = "  +4.59%  " }
    @{ Row = 13; D = $null;        E = "  +1.49%  " }
    @{ Row = 14; D = "3.657.20";   E = "  +0.42%  " }
    @{ Row = 15; D = "25.98";      E = "  +2.64%  " }
    @{ Row = 16; D = $null;        E = "  +0.68%  " }
    @{ Row = 17; D = "57.797.99";  E = "  +0.37%  " }
    @{ Row = 18; D = "3.124.98";   E = "  +0.75%  " }
    @{ Row = 19; D = "6.06";       E = "  +2.15%  " }
    @{ Row = 20; D = "12.70";      E = "  +2.70%  " }
    @{ Row = 21; D = "8.07";       E = "  +2.59%  " }
    @{ Row = 22; D = "366.31";     E = "  +6.17%  " }
    @{ Row = 23; D = "0.998";      E = "  -0.38%  " }
    @{ Row = 24; D = "5.68";       E = "  -1.80%  " }
    @{ Row = 25; D = "69.10";      E = "  +2.10%  " }
    @{ Row = 26; D = "0.505";      E = "  +1.03%  " }
    @{ Row = 27; D = $null;        E = "  +0.34%  " }
    @{ Row = 28; D = $null;        E = "  +0.26%  " }
    @{ Row = 29; D = "0.0₃0861";   E = "  -3.37%  " }
    @{ Row = 30; D = $null;        E = "  -1.75%  " }
    @{ Row = 31; D = $null;        E = "  -0.11%  " }
    @{ Row = 32; D = "6.05";       E = "  +0.37%  " }
    @{ Row = 33; D = "21.35";      E = "  +2.64%  " }
    @{ Row = 34; D = $null;        E = "  +3.79%  " }
    @{ Row = 35; D = $null;        E = "  -0.50%  " }
    @{ Row = 36; D = "159.49";     E = "  +0.69%  " }
    @{ Row = 37; D = $null;        E = "  +0.09%  " }
    @{ Row = 38; D = $null;        E = "  +5.21%  " }
    @{ Row = 39; D = "25.40";      E = "  -2.12%  " }
    @{ Row = 40; D = "1.67";       E = "  +3.38%  " }
    @{ Row = 41; D = "0.0671";     E = "  +1.63%  " }
    @{ Row = 42; D = "2.536.32";   E = "  +6.86%  " }
    @{ Row = 43; D = "4.04";       E = "  -1.31%  " }
    @{ Row = 44; D = "0.698";      E = "  -0.15%  " }
    @{ Row = 45; D = $null;        E = "  +3.22%  " }
    @{ Row = 47; D = $null;        E = "  +0.56%  " }
    @{ Row = 48; D = "0.980";      E = "  +0.61%  " }
    @{ Row = 49; D = $null;        E = "  +2.06%  " }
    @{ Row = 50; D = "19.68";      E = "  -0.47%  " }
    @{ Row = 51; D = "0.738";      E = "  -2.78%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        # Force the cell to remain plain text (the price strings use "."
        # as a thousands separator, e.g. "57.691.84", and some plain
        # decimals like "12.70" would otherwise be auto-coerced into a
        # Number and lose the trailing zero). Reset the style back to
        # Normal afterwards so no stray number-format is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
